# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.017.07"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.587.18"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'522.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'139.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.58%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.597.46"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "3.045.39"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "58.949.92"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'20.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "2.582.98"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "'340.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "'4.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").Value = "'6.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'66.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'7.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0722"
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("E31").Value = "  -6.13%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "'18.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "'149.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'3.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "'1.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").Value = "'36.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D39").Value = "'0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.39%  "
$ws.Range("D41").Value = "'3.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'272.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'10.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "'0.595"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'0.0950"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'0.0516"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "'18.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").Value = "1.970.30"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -4.81%  "
